# Add a trailing "description" column to the header row of every worksheet
# that models a LinkML "AttributeGroup" (a complex-attribute value object
# without its own independent identity) rather than a "NamedThing" (which
# already carries id/title/description columns and is left untouched).
#
# Sheets touched (header row 1, new column appended right after the last
# existing column):
#   ImageFeature, MolecularComposition, BufferComposition,
#   StorageConditions, CryoEMPreparation, XRayPreparation,
#   SAXSPreparation, ExperimentalConditions, DataCollectionStrategy,
#   QualityMetrics, ComputeResources

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "ImageFeature",
    "MolecularComposition",
    "BufferComposition",
    "StorageConditions",
    "CryoEMPreparation",
    "XRayPreparation",
    "SAXSPreparation",
    "ExperimentalConditions",
    "DataCollectionStrategy",
    "QualityMetrics",
    "ComputeResources"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Number of already-populated header columns in row 1.
    $lastColumnIndex = $ws.UsedRange.Columns.Count

    $newColumnIndex = $lastColumnIndex + 1
    $ws.Cells.Item(1, $newColumnIndex).Value2 = "description"
}

Write-Host "Appended 'description' column to $($sheetNames.Count) sheets"
